$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("D2").Value = 1.8099369330201
$ws.Range("F2").Value = 0.1842

# Update row 3
$ws.Range("D3").Value = 2.34335772134835
$ws.Range("F3").Value = 0.129

# Row 4 gets the values that used to belong to row 5, with updated D/E/F
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0.14527443557455
$ws.Range("C4").Value = 0.14527443557455
$ws.Range("D4").Value = 0.132124429783464
$ws.Range("E4").Value = 0.000891014306055836
$ws.Range("F4").Value = 0.7176

# Row 5 becomes new data (previously row "132" entries), D and F now blank
$ws.Range("A5").Value = 144
$ws.Range("B5").Value = 158.331950851328
$ws.Range("C5").Value = 1.09952743646756
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 0.97110019912531
$ws.Range("F5").ClearContents()

# Row 6 becomes the previous row 147 entry, C/D/F blank
$ws.Range("A6").Value = 147
$ws.Range("B6").Value = 163.043886711115
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 1
$ws.Range("F6").ClearContents()

# Remove old rows 7-10 entirely (data no longer present)
$ws.Range("A7:F10").ClearContents()
